# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing data, mirroring the header style already used in AC1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, bordered, centered) from AC1
# onto the three new header cells so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row shares the same team season record.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 64
    $ws.Cells.Item($row, 31).Value = 98
    $ws.Cells.Item($row, 32).Value = 0
}
